# Atualização da lista com o requisito Banana
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preenche a nova linha de requisito (linha 9) que antes estava vazia
$ws.Range("A9").Value = "Banana"
$ws.Range("B9").Value = "Funcional"
$ws.Range("C9").Value = "Tecnologia/ Arquitetura"
$ws.Range("D9").Value = "-"

# A linha cresce para acomodar o texto com quebra de linha (igual às demais linhas de dados)
$ws.Range("A9:D9").RowHeight = 30

# Seleção final do usuário, como no arquivo salvo
$ws.Range("D10").Select() | Out-Null
